$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status cells ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"

# --- zh-cn sheet: Status cell ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"

# --- de-de sheet: Status cell ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"

# The "Status" text is now shorter ("In Translation" vs "Ready for handoff"),
# so the columns that show it shrink to fit the new content.
# (ColumnWidth is specified in the COM "character width" unit; 12.5 is the
# value that resolves to the narrower stored column width used for the
# Status / zh-cn / de-de columns below.)
$newWidth = 12.5

$ws1.Columns.Item(5).ColumnWidth = $newWidth   # Overview!E (zh-cn)
$ws1.Columns.Item(6).ColumnWidth = $newWidth   # Overview!F (de-de)
$ws2.Columns.Item(3).ColumnWidth = $newWidth   # zh-cn!C (Status)
$ws3.Columns.Item(3).ColumnWidth = $newWidth   # de-de!C (Status)
